$d = $word.ActiveDocument

# Locate the "Places" text (table cell) that needs to become "Place".
$r = $d.Content
$r.Find.Execute("Places", $true, $false, $false, $false, $false,
                $true, 1, $false, "", 0)
$placesStart = $r.Start
$placesEnd = $r.End

# Word automatically (re)positions its hidden "_GoBack" bookmark at the
# location of the most recent edit. Recreate/move it here: a zero-length
# bookmark sitting right after "Place" (i.e. right before the trailing
# "s" that is about to be removed). Adding a bookmark with a name that is
# already in use elsewhere in the document relocates the existing one.
$bmPos = $placesEnd - 1
$bm = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bm)

# Turn "Places" into "Place" by removing the trailing "s".
$trailing = $d.Range($placesEnd - 1, $placesEnd)
$trailing.Delete()
